$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$day15 = "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/Day 15/"
$day16 = "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/Day 16/"

# --- Day labels (A16, A17) first, matching shared-string creation order ---
$ws.Cells.Item(16,1).Value = "Day 15"
$ws.Cells.Item(17,1).Value = "Day 16"

# --- Dates (B16, B17); copy date format from B14 afterward so the existing
#     date style (s=2) is reused instead of a brand-new numFmt being minted ---
$ws.Cells.Item(16,2).Value = 45817
$ws.Cells.Item(17,2).Value = 45818
$ws.Cells.Item(14,2).Copy() | Out-Null
$ws.Cells.Item(16,2).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(17,2).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Row 16 (Day 15) hyperlinked question cells ---
$ws.Hyperlinks.Add($ws.Cells.Item(16,3), ($day15 + "Isomorphic Strings.py"), "", "Isomorphic Strings.py", ($day15 + "Isomorphic Strings.py")) | Out-Null
$ws.Cells.Item(16,3).Value = "Isomorphic Strings.py"

$ws.Hyperlinks.Add($ws.Cells.Item(16,4), ($day15 + "Simplify Path.py"), "", "Simplify Path.py", ($day15 + "Simplify Path.py")) | Out-Null
$ws.Cells.Item(16,4).Value = "Simplify Path.py"

$ws.Hyperlinks.Add($ws.Cells.Item(16,5), ($day15 + "Time Based Key-Value Store.py"), "", "Time Based Key-Value Store.py", ($day15 + "Time Based Key-Value Store.py")) | Out-Null
$ws.Cells.Item(16,5).Value = "Time Based Key-Value Store.py"

# --- Row 17 (Day 16) hyperlinked question cells ---
$ws.Hyperlinks.Add($ws.Cells.Item(17,3), ($day16 + "Insert Delete GetRandom O(1).py"), "", "Insert Delete GetRandom O(1).py", ($day16 + "Insert Delete GetRandom O(1).py")) | Out-Null
$ws.Cells.Item(17,3).Value = "Insert Delete GetRandom O(1).py"

$ws.Hyperlinks.Add($ws.Cells.Item(17,4), ($day16 + "Reverse Words in a String.py"), "", "Reverse Words in a String.py", ($day16 + "Reverse Words in a String.py")) | Out-Null
$ws.Cells.Item(17,4).Value = "Reverse Words in a String.py"

$ws.Hyperlinks.Add($ws.Cells.Item(17,5), ($day16 + "Zigzag Conversion.py"), "", "Zigzag Conversion.py", ($day16 + "Zigzag Conversion.py")) | Out-Null
$ws.Cells.Item(17,5).Value = "Zigzag Conversion.py"

# Re-apply the existing hyperlink cell style (font/border/etc from C14:E14)
# onto the new cells so no duplicate style entries are created.
$ws.Cells.Item(14,3).Copy() | Out-Null
$ws.Range("C16:E17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Topics (F16, F17) and status (G16, G17) last ---
$ws.Cells.Item(16,6).Value = "Stack, Binary Search, Stack"
$ws.Cells.Item(17,6).Value = "String/2Pointer, Stack, HashMap,Arrays"
$ws.Cells.Item(16,7).Value = "S"
$ws.Cells.Item(17,7).Value = "S"

$excel.CutCopyMode = $false

$ws.Range("H14").Select() | Out-Null

Write-Output "done"
